$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 25000168
$ws.Cells.Item(4, 9).Value = 191.85715
$ws.Cells.Item(4, 10).Value = 83333450
$ws.Cells.Item(4, 11).Value = 191.85715
$ws.Cells.Item(4, 12).Value = 83333450
$ws.Cells.Item(4, 13).Value = -77.85714999999999
$ws.Cells.Item(4, 14).Value = -83333678
$ws.Cells.Item(74, 8).Value = 4142.143
$ws.Cells.Item(74, 9).Value = 2997.5
$ws.Cells.Item(74, 11).Value = 2997.5
$ws.Cells.Item(74, 13).Value = -2061.5
$ws.Cells.Item(77, 8).Value = 4142.143
$ws.Cells.Item(77, 9).Value = 2997.5
$ws.Cells.Item(77, 11).Value = 14987.5
$ws.Cells.Item(77, 13).Value = -10307.5
$ws.Cells.Item(94, 8).Value = 2276.5
$ws.Cells.Item(94, 9).Value = 1173.1428
$ws.Cells.Item(94, 10).Value = 10000
$ws.Cells.Item(94, 11).Value = 1173.1428
$ws.Cells.Item(94, 12).Value = 10000
$ws.Cells.Item(94, 13).Value = -722.1428000000001
$ws.Cells.Item(94, 14).Value = -10902
$ws.Cells.Item(106, 8).Value = 5394
$ws.Cells.Item(106, 9).Value = 1521.7142
$ws.Cells.Item(106, 11).Value = 1521.7142
$ws.Cells.Item(106, 13).Value = -890.7141999999999
$ws.Cells.Item(112, 8).Value = 1956.4857
$ws.Cells.Item(112, 9).Value = 7000
$ws.Cells.Item(112, 11).Value = 21000
$ws.Cells.Item(112, 13).Value = -19892
$ws.Cells.Item(113, 8).Value = 76926810
$ws.Cells.Item(113, 9).Value = 200002610
$ws.Cells.Item(113, 10).Value = 4436.375
$ws.Cells.Item(113, 11).Value = 200002610
$ws.Cells.Item(113, 12).Value = 4436.375
$ws.Cells.Item(113, 13).Value = -199999356
$ws.Cells.Item(113, 14).Value = -10944.375
$ws.Cells.Item(138, 8).Value = 2458.7073
$ws.Cells.Item(138, 10).Value = 2799.2068
$ws.Cells.Item(138, 12).Value = 8397.6204
$ws.Cells.Item(138, 14).Value = -18677.6204

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1764.9493
$ws.Cells.Item(32, 9).Value = 1715.4868
$ws.Cells.Item(32, 11).Value = 1715.4868
$ws.Cells.Item(32, 13).Value = -1428.4868
$ws.Cells.Item(45, 8).Value = 2263.6667
$ws.Cells.Item(45, 9).Value = 1924.6
$ws.Cells.Item(45, 11).Value = 1924.6
$ws.Cells.Item(45, 13).Value = -1547.6
$ws.Cells.Item(61, 8).Value = 16668561
$ws.Cells.Item(61, 9).Value = 20834784
$ws.Cells.Item(61, 10).Value = 3666.25
$ws.Cells.Item(61, 11).Value = 20834784
$ws.Cells.Item(61, 12).Value = 3666.25
$ws.Cells.Item(61, 13).Value = -20834572
$ws.Cells.Item(61, 14).Value = -4090.25
$ws.Cells.Item(88, 8).Value = 6174174
$ws.Cells.Item(88, 9).Value = 13889431
$ws.Cells.Item(88, 10).Value = 1968.4667
$ws.Cells.Item(88, 11).Value = 13889431
$ws.Cells.Item(88, 12).Value = 1968.4667
$ws.Cells.Item(88, 13).Value = -13889025
$ws.Cells.Item(88, 14).Value = -2780.4667
$ws.Cells.Item(91, 8).Value = 6174174
$ws.Cells.Item(91, 9).Value = 13889431
$ws.Cells.Item(91, 10).Value = 1968.4667
$ws.Cells.Item(91, 11).Value = 13889431
$ws.Cells.Item(91, 12).Value = 1968.4667
$ws.Cells.Item(91, 13).Value = -13888027
$ws.Cells.Item(91, 14).Value = -4776.4667
$ws.Cells.Item(110, 8).Value = 66687864
$ws.Cells.Item(110, 10).Value = 75199.75
$ws.Cells.Item(110, 12).Value = 75199.75
$ws.Cells.Item(110, 14).Value = -79289.75
$ws.Cells.Item(136, 8).Value = 16668561
$ws.Cells.Item(136, 9).Value = 20834784
$ws.Cells.Item(136, 10).Value = 3666.25
$ws.Cells.Item(136, 11).Value = 62504352
$ws.Cells.Item(136, 12).Value = 10998.75
$ws.Cells.Item(136, 13).Value = -62501802
$ws.Cells.Item(136, 14).Value = -16098.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 14706973
$ws.Cells.Item(107, 9).Value = 1162.4286
$ws.Cells.Item(107, 10).Value = 83334090
$ws.Cells.Item(107, 11).Value = 1162.4286
$ws.Cells.Item(107, 12).Value = 83334090
$ws.Cells.Item(107, 13).Value = 757.5714
$ws.Cells.Item(107, 14).Value = -83337930
$ws.Cells.Item(134, 8).Value = 1767.7858
$ws.Cells.Item(134, 9).Value = 1345.75
$ws.Cells.Item(134, 10).Value = 4300
$ws.Cells.Item(134, 11).Value = 4037.25
$ws.Cells.Item(134, 12).Value = 12900
$ws.Cells.Item(134, 13).Value = -1502.25
$ws.Cells.Item(134, 14).Value = -17970

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1578.8636
$ws.Cells.Item(16, 9).Value = 1406
$ws.Cells.Item(16, 10).Value = 1786.3
$ws.Cells.Item(16, 11).Value = 1406
$ws.Cells.Item(16, 12).Value = 1786.3
$ws.Cells.Item(16, 13).Value = -1119
$ws.Cells.Item(16, 14).Value = -2360.3
$ws.Cells.Item(22, 8).Value = 172.2
$ws.Cells.Item(22, 9).Value = 87
$ws.Cells.Item(22, 10).Value = 300
$ws.Cells.Item(22, 11).Value = 87
$ws.Cells.Item(22, 12).Value = 300
$ws.Cells.Item(22, 13).Value = 263
$ws.Cells.Item(22, 14).Value = -1000
$ws.Cells.Item(86, 8).Value = 45796
$ws.Cells.Item(86, 9).Value = 59614.668
$ws.Cells.Item(86, 10).Value = 4340
$ws.Cells.Item(86, 11).Value = 59614.668
$ws.Cells.Item(86, 12).Value = 4340
$ws.Cells.Item(86, 13).Value = -58491.668
$ws.Cells.Item(86, 14).Value = -6586
$ws.Cells.Item(89, 8).Value = 45796
$ws.Cells.Item(89, 9).Value = 59614.668
$ws.Cells.Item(89, 10).Value = 4340
$ws.Cells.Item(89, 11).Value = 298073.34
$ws.Cells.Item(89, 12).Value = 21700
$ws.Cells.Item(89, 13).Value = -292457.34
$ws.Cells.Item(89, 14).Value = -32932
$ws.Cells.Item(113, 8).Value = 1578.8636
$ws.Cells.Item(113, 9).Value = 1406
$ws.Cells.Item(113, 10).Value = 1786.3
$ws.Cells.Item(113, 11).Value = 1406
$ws.Cells.Item(113, 12).Value = 1786.3
$ws.Cells.Item(113, 13).Value = 764
$ws.Cells.Item(113, 14).Value = -6126.3
$ws.Cells.Item(122, 8).Value = 1414.6666
$ws.Cells.Item(122, 9).Value = 1414.6666
$ws.Cells.Item(122, 11).Value = 4243.9998
$ws.Cells.Item(122, 13).Value = -1793.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10320.88
$ws.Cells.Item(70, 9).Value = 9768.333000000001
$ws.Cells.Item(70, 11).Value = 9768.333000000001
$ws.Cells.Item(70, 13).Value = -9498.333000000001
$ws.Cells.Item(73, 8).Value = 10320.88
$ws.Cells.Item(73, 9).Value = 9768.333000000001
$ws.Cells.Item(73, 11).Value = 9768.333000000001
$ws.Cells.Item(73, 13).Value = -8832.333000000001
$ws.Cells.Item(97, 8).Value = 2070.3914
$ws.Cells.Item(97, 9).Value = 2015.2354
$ws.Cells.Item(97, 10).Value = 2226.6667
$ws.Cells.Item(97, 11).Value = 2015.2354
$ws.Cells.Item(97, 12).Value = 2226.6667
$ws.Cells.Item(97, 13).Value = -1519.2354
$ws.Cells.Item(97, 14).Value = -3218.6667
$ws.Cells.Item(113, 8).Value = 3619.2222
$ws.Cells.Item(113, 10).Value = 4287.5713
$ws.Cells.Item(113, 12).Value = 4287.5713
$ws.Cells.Item(113, 14).Value = -8627.5713
$ws.Cells.Item(122, 8).Value = 45458016
$ws.Cells.Item(122, 9).Value = 2047.25
$ws.Cells.Item(122, 10).Value = 71432856
$ws.Cells.Item(122, 11).Value = 6141.75
$ws.Cells.Item(122, 12).Value = 214298568
$ws.Cells.Item(122, 13).Value = -3691.75
$ws.Cells.Item(122, 14).Value = -214303468
$ws.Cells.Item(126, 8).Value = 9436.235000000001
$ws.Cells.Item(126, 9).Value = 19168.666
$ws.Cells.Item(126, 10).Value = 4127.636
$ws.Cells.Item(126, 11).Value = 57505.99800000001
$ws.Cells.Item(126, 12).Value = 12382.908
$ws.Cells.Item(126, 13).Value = -55035.99800000001
$ws.Cells.Item(126, 14).Value = -17322.908
$ws.Cells.Item(127, 8).Value = 326
$ws.Cells.Item(127, 10).Value = 326
$ws.Cells.Item(127, 12).Value = 326
$ws.Cells.Item(127, 14).Value = -10246

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 31253768
$ws.Cells.Item(7, 9).Value = 55558430
$ws.Cells.Item(7, 11).Value = 55558430
$ws.Cells.Item(7, 13).Value = -55558318
$ws.Cells.Item(22, 8).Value = 3226195
$ws.Cells.Item(22, 9).Value = 300
$ws.Cells.Item(22, 10).Value = 4608721.5
$ws.Cells.Item(22, 11).Value = 300
$ws.Cells.Item(22, 12).Value = 4608721.5
$ws.Cells.Item(22, 13).Value = -5
$ws.Cells.Item(22, 14).Value = -4609311.5
$ws.Cells.Item(27, 8).Value = 3226195
$ws.Cells.Item(27, 9).Value = 300
$ws.Cells.Item(27, 10).Value = 4608721.5
$ws.Cells.Item(27, 11).Value = 300
$ws.Cells.Item(27, 12).Value = 4608721.5
$ws.Cells.Item(27, 13).Value = -193
$ws.Cells.Item(27, 14).Value = -4608935.5
$ws.Cells.Item(31, 8).Value = 1160.2307
$ws.Cells.Item(31, 9).Value = 799.3333
$ws.Cells.Item(31, 10).Value = 1972.25
$ws.Cells.Item(31, 11).Value = 799.3333
$ws.Cells.Item(31, 12).Value = 1972.25
$ws.Cells.Item(31, 13).Value = -551.3333
$ws.Cells.Item(31, 14).Value = -2468.25
$ws.Cells.Item(61, 8).Value = 4986.385
$ws.Cells.Item(61, 9).Value = 3387.6
$ws.Cells.Item(61, 10).Value = 10315.667
$ws.Cells.Item(61, 11).Value = 3387.6
$ws.Cells.Item(61, 12).Value = 10315.667
$ws.Cells.Item(61, 13).Value = -3185.6
$ws.Cells.Item(61, 14).Value = -10719.667
$ws.Cells.Item(68, 8).Value = 6107.92
$ws.Cells.Item(68, 9).Value = 2446.2307
$ws.Cells.Item(68, 10).Value = 10074.75
$ws.Cells.Item(68, 11).Value = 2446.2307
$ws.Cells.Item(68, 12).Value = 10074.75
$ws.Cells.Item(68, 13).Value = -1697.2307
$ws.Cells.Item(68, 14).Value = -11572.75
$ws.Cells.Item(71, 8).Value = 6107.92
$ws.Cells.Item(71, 9).Value = 2446.2307
$ws.Cells.Item(71, 10).Value = 10074.75
$ws.Cells.Item(71, 11).Value = 12231.1535
$ws.Cells.Item(71, 12).Value = 50373.75
$ws.Cells.Item(71, 13).Value = -8487.1535
$ws.Cells.Item(71, 14).Value = -57861.75
$ws.Cells.Item(93, 8).Value = 3333.476
$ws.Cells.Item(93, 9).Value = 3258
$ws.Cells.Item(93, 10).Value = 3522.1667
$ws.Cells.Item(93, 11).Value = 3258
$ws.Cells.Item(93, 12).Value = 3522.1667
$ws.Cells.Item(93, 13).Value = -2010
$ws.Cells.Item(93, 14).Value = -6018.1667
$ws.Cells.Item(113, 8).Value = 4986.385
$ws.Cells.Item(113, 9).Value = 3387.6
$ws.Cells.Item(113, 10).Value = 10315.667
$ws.Cells.Item(113, 11).Value = 3387.6
$ws.Cells.Item(113, 12).Value = 10315.667
$ws.Cells.Item(113, 13).Value = -1217.6
$ws.Cells.Item(113, 14).Value = -14655.667
$ws.Cells.Item(122, 8).Value = 3686.2593
$ws.Cells.Item(122, 9).Value = 2908.1875
$ws.Cells.Item(122, 10).Value = 4818
$ws.Cells.Item(122, 11).Value = 8724.5625
$ws.Cells.Item(122, 12).Value = 14454
$ws.Cells.Item(122, 13).Value = -6274.5625
$ws.Cells.Item(122, 14).Value = -19354
$ws.Cells.Item(126, 8).Value = 31253768
$ws.Cells.Item(126, 9).Value = 55558430
$ws.Cells.Item(126, 11).Value = 166675290
$ws.Cells.Item(126, 13).Value = -166672820
$ws.Cells.Item(136, 8).Value = 2399.4285
$ws.Cells.Item(136, 10).Value = 3023
$ws.Cells.Item(136, 12).Value = 9069
$ws.Cells.Item(136, 14).Value = -14169

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 536.1875
$ws.Cells.Item(107, 9).Value = 555.26666
$ws.Cells.Item(107, 10).Value = 250
$ws.Cells.Item(107, 11).Value = 1665.79998
$ws.Cells.Item(107, 12).Value = 750
$ws.Cells.Item(107, 13).Value = 254.20002
$ws.Cells.Item(107, 14).Value = -4590
$ws.Cells.Item(113, 8).Value = 10733.167
$ws.Cells.Item(113, 9).Value = 200
$ws.Cells.Item(113, 10).Value = 15999.75
$ws.Cells.Item(113, 11).Value = 600
$ws.Cells.Item(113, 12).Value = 47999.25
$ws.Cells.Item(113, 13).Value = 1570
$ws.Cells.Item(113, 14).Value = -52339.25
$ws.Cells.Item(132, 8).Value = 7500.25
$ws.Cells.Item(132, 9).Value = 9167.333000000001
$ws.Cells.Item(132, 10).Value = 5833.1665
$ws.Cells.Item(132, 11).Value = 27501.999
$ws.Cells.Item(132, 12).Value = 17499.4995
$ws.Cells.Item(132, 13).Value = -24971.999
$ws.Cells.Item(132, 14).Value = -22559.4995
$ws.Cells.Item(136, 8).Value = 8249.833000000001
$ws.Cells.Item(136, 9).Value = 3000
$ws.Cells.Item(136, 11).Value = 9000
$ws.Cells.Item(136, 13).Value = -6450
